# Export support option done data to excel file.
# Add a new "Support Option Done" worksheet positioned right before the
# existing "ART Refill" sheet, with the sheet title written to cell A1
# (matching the single-cell "title" layout used by the other sheets in
# this workbook).

$wb = $excel.ActiveWorkbook

$artRefill = $wb.Worksheets.Item("ART Refill")
$ws = $wb.Worksheets.Add($artRefill)
$ws.Name = "Support Option Done"
$ws.Range("A1").Value = "Support Option Done"

# Adding a sheet makes it the active one; restore the originally active
# sheet so the rest of the workbook's view state is left untouched.
$wb.Worksheets.Item("User Data").Activate()
